$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 6258794
$ws.Range("J17").Value = 6675920.5
$ws.Range("L17").Value = 20027761.5
$ws.Range("N17").Value = -20028097.5

# Row 88
$ws.Range("H88").Value = 661.8125
$ws.Range("I88").Value = 317.66666
$ws.Range("J88").Value = 741.2308
$ws.Range("K88").Value = 317.66666
$ws.Range("L88").Value = 741.2308
$ws.Range("M88").Value = 88.33334000000002
$ws.Range("N88").Value = -1553.2308

# Row 91
$ws.Range("H91").Value = 661.8125
$ws.Range("I91").Value = 317.66666
$ws.Range("J91").Value = 741.2308
$ws.Range("K91").Value = 317.66666
$ws.Range("L91").Value = 741.2308
$ws.Range("M91").Value = 1086.33334
$ws.Range("N91").Value = -3549.2308

# Row 103
$ws.Range("H103").Value = 399.5
$ws.Range("I103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("M103").ClearContents()

# Row 106
$ws.Range("H106").Value = 10755773
$ws.Range("I106").Value = 41668180
$ws.Range("J106").Value = 3631.739
$ws.Range("K106").Value = 41668180
$ws.Range("L106").Value = 3631.739
$ws.Range("M106").Value = -41667549
$ws.Range("N106").Value = -4893.739

# Row 107
$ws.Range("H107").Value = 697.24
$ws.Range("I107").Value = 583.95654
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 583.95654
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = 1336.04346
$ws.Range("N107").Value = -5840

# Row 129
$ws.Range("H129").Value = 141767.02
$ws.Range("J129").Value = 167713.2
$ws.Range("L129").Value = 503139.6
$ws.Range("N129").Value = -513139.6

# Row 131
$ws.Range("H131").Value = 1694.1428
$ws.Range("J131").Value = 2200
$ws.Range("L131").Value = 6600
$ws.Range("N131").Value = -16680

# Row 132
$ws.Range("H132").Value = 3585.7693
$ws.Range("I132").Value = 3879.5652
$ws.Range("K132").Value = 11638.6956
$ws.Range("M132").Value = -9108.695599999999

# Row 141
$ws.Range("H141").Value = 3282.923
$ws.Range("I141").Value = 3129
$ws.Range("J141").Value = 3414.8572
$ws.Range("K141").Value = 9387
$ws.Range("L141").Value = 10244.5716
$ws.Range("M141").Value = -4207
$ws.Range("N141").Value = -20604.5716


$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 735.325
$ws.Range("I2").Value = 666.25
$ws.Range("K2").Value = 666.25
$ws.Range("M2").Value = -553.25

# Row 32
$ws.Range("H32").Value = 5470.3374
$ws.Range("I32").Value = 4161.2056
$ws.Range("J32").Value = 12821.615
$ws.Range("K32").Value = 4161.2056
$ws.Range("L32").Value = 12821.615
$ws.Range("M32").Value = -3874.2056
$ws.Range("N32").Value = -13395.615

# Row 116
$ws.Range("H116").Value = 735.325
$ws.Range("I116").Value = 666.25
$ws.Range("K116").Value = 666.25
$ws.Range("M116").Value = 1627.75

# Row 122
$ws.Range("H122").Value = 1790.15
$ws.Range("I122").Value = 1610.6842
$ws.Range("K122").Value = 4832.0526
$ws.Range("M122").Value = -2382.0526

# Row 132
$ws.Range("H132").Value = 10501.842
$ws.Range("I132").Value = 1641.2222
$ws.Range("J132").Value = 25691.477
$ws.Range("K132").Value = 4923.6666
$ws.Range("L132").Value = 77074.431
$ws.Range("M132").Value = -2393.6666
$ws.Range("N132").Value = -82134.431


$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 735.325
$ws.Range("I3").Value = 666.25
$ws.Range("K3").Value = 666.25
$ws.Range("M3").Value = -552.25

# Row 94
$ws.Range("H94").Value = 691.8570999999999
$ws.Range("I94").Value = 454.15384
$ws.Range("J94").Value = 1078.125
$ws.Range("K94").Value = 454.15384
$ws.Range("L94").Value = 1078.125
$ws.Range("M94").Value = -3.153840000000002
$ws.Range("N94").Value = -1980.125

# Row 105
$ws.Range("H105").Value = 2501897.5
$ws.Range("I105").Value = 1812.5
$ws.Range("J105").Value = 3126918.8
$ws.Range("K105").Value = 1812.5
$ws.Range("L105").Value = 3126918.8
$ws.Range("M105").Value = -65.5
$ws.Range("N105").Value = -3130412.8

# Row 134
$ws.Range("H134").Value = 4190.6895
$ws.Range("I134").Value = 4308.2144
$ws.Range("J134").Value = 900
$ws.Range("K134").Value = 12924.6432
$ws.Range("L134").Value = 2700
$ws.Range("M134").Value = -7770


$ws = $wb.Worksheets.Item("CRP")
# Row 107
$ws.Range("H107").Value = 1815.8334
$ws.Range("I107").Value = 617.25
$ws.Range("J107").Value = 2774.7
$ws.Range("K107").Value = 617.25
$ws.Range("L107").Value = 2774.7
$ws.Range("M107").Value = 1302.75
$ws.Range("N107").Value = -6614.7

# Row 140
$ws.Range("H140").Value = 41550
$ws.Range("J140").Value = 41550
$ws.Range("L140").Value = 41550
$ws.Range("N140").Value = -51910


$ws = $wb.Worksheets.Item("CUL")
# Row 21
$ws.Range("H21").Value = 524.6
$ws.Range("I21").Value = 311.5
$ws.Range("J21").Value = 666.6667
$ws.Range("K21").Value = 934.5
$ws.Range("L21").Value = 2000.0001
$ws.Range("M21").Value = -2346.0001

# Row 39
$ws.Range("H39").Value = 3373.1
$ws.Range("J39").Value = 3373.1
$ws.Range("L39").Value = 10119.3
$ws.Range("N39").Value = -10707.3

# Row 56
$ws.Range("H56").Value = 6816.154
$ws.Range("I56").Value = 6816.154
$ws.Range("K56").Value = 6816.154
$ws.Range("M56").Value = -6286.154

# Row 131
$ws.Range("H131").Value = 708.8099999999999
$ws.Range("J131").Value = 723.20215
$ws.Range("L131").Value = 2169.60645
$ws.Range("N131").Value = -12249.60645

# Row 133
$ws.Range("H133").Value = 5250
$ws.Range("I133").Value = 1377.5
$ws.Range("J133").Value = 6799
$ws.Range("K133").Value = 4132.5
$ws.Range("L133").Value = 20397
$ws.Range("M133").Value = 927.5
$ws.Range("N133").Value = -30517

# Row 141
$ws.Range("H141").Value = 5532.6665
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()


$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3721.0833
$ws.Range("J80").Value = 4007.5715
$ws.Range("L80").Value = 4007.5715
$ws.Range("N80").Value = -6003.5715

# Row 83
$ws.Range("H83").Value = 3721.0833
$ws.Range("J83").Value = 4007.5715
$ws.Range("L83").Value = 20037.8575
$ws.Range("N83").Value = -30021.8575

# Row 113
$ws.Range("H113").Value = 10522.818
$ws.Range("I113").Value = 12968.875
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 12968.875
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = -10798.875
$ws.Range("N113").Value = -8340

# Row 132
$ws.Range("H132").Value = 19590.406
$ws.Range("I132").Value = 4310.8096
$ws.Range("J132").Value = 48760.547
$ws.Range("K132").Value = 12932.4288
$ws.Range("L132").Value = 146281.641
$ws.Range("M132").Value = -10402.4288
$ws.Range("N132").Value = -151341.641


$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 559.1111
$ws.Range("J16").Value = 887.75
$ws.Range("L16").Value = 887.75
$ws.Range("N16").Value = -1227.75

# Row 36
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").ClearContents()
$ws.Range("N36").Value = 0

# Row 40
$ws.Range("H40").Value = 2809.8823
$ws.Range("I40").Value = 2579.6667
$ws.Range("J40").Value = 3697.8572
$ws.Range("K40").Value = 2579.6667
$ws.Range("L40").Value = 3697.8572
$ws.Range("M40").Value = -2443.6667
$ws.Range("N40").Value = -3969.8572

# Row 46
$ws.Range("H46").Value = 2057.48
$ws.Range("I46").Value = 2338.2144
$ws.Range("J46").Value = 1700.1818
$ws.Range("K46").Value = 2338.2144
$ws.Range("L46").Value = 1700.1818
$ws.Range("M46").Value = -2150.2144
$ws.Range("N46").Value = -2076.1818

# Row 93
$ws.Range("H93").Value = 1446.0769
$ws.Range("I93").Value = 1389.9
$ws.Range("J93").Value = 1633.3334
$ws.Range("K93").Value = 1389.9
$ws.Range("L93").Value = 1633.3334
$ws.Range("M93").Value = -141.9000000000001
$ws.Range("N93").Value = -4129.3334


$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 954.4286
$ws.Range("I122").Value = 699.86664
$ws.Range("J122").Value = 1590.8334
$ws.Range("K122").Value = 2099.59992
$ws.Range("L122").Value = 4772.5002
$ws.Range("M122").Value = 350.4000800000003
$ws.Range("N122").Value = -9672.5002

# Row 132
$ws.Range("H132").Value = 1328.1282
$ws.Range("I132").Value = 872.5484
$ws.Range("J132").Value = 3093.5
$ws.Range("K132").Value = 2617.6452
$ws.Range("L132").Value = 9280.5
$ws.Range("M132").Value = -87.64519999999993
$ws.Range("N132").Value = -14340.5

# Row 136
$ws.Range("H136").Value = 25180072
$ws.Range("I136").Value = 32259452
$ws.Range("J136").Value = 8945
$ws.Range("K136").Value = 96778356
$ws.Range("L136").Value = 26835
$ws.Range("M136").Value = -96775806
$ws.Range("N136").Value = -31935

